$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 100815.664
$ws.Range("J133").Value = 100815.664
$ws.Range("L133").Value = 100815.664
$ws.Range("N133").Value = -110935.664
$ws.Range("H138").Value = 3340.7966
$ws.Range("J138").Value = 4003.375
$ws.Range("L138").Value = 12010.125
$ws.Range("N138").Value = -22290.125

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1328.9333
$ws.Range("I2").Value = 1066.7142
$ws.Range("K2").Value = 1066.7142
$ws.Range("M2").Value = -953.7141999999999
$ws.Range("H5").Value = 231.71428
$ws.Range("I5").Value = 278
$ws.Range("K5").Value = 278
$ws.Range("M5").Value = -166
$ws.Range("H32").Value = 24327.916
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 24327.916
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 24327.916
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -24901.916
$ws.Range("H61").Value = 3101.1904
$ws.Range("I61").Value = 1700
$ws.Range("K61").Value = 1700
$ws.Range("M61").Value = -1488
$ws.Range("H94").Value = 6894.5
$ws.Range("J94").Value = 6894.5
$ws.Range("L94").Value = 6894.5
$ws.Range("N94").Value = -8696.5
$ws.Range("H116").Value = 1328.9333
$ws.Range("I116").Value = 1066.7142
$ws.Range("K116").Value = 1066.7142
$ws.Range("M116").Value = 1227.2858
$ws.Range("H122").Value = 1950.5454
$ws.Range("I122").Value = 1950.5454
$ws.Range("K122").Value = 5851.6362
$ws.Range("M122").Value = -3401.6362
$ws.Range("H136").Value = 3101.1904
$ws.Range("I136").Value = 1700
$ws.Range("K136").Value = 5100
$ws.Range("M136").Value = -2550

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1328.9333
$ws.Range("I3").Value = 1066.7142
$ws.Range("K3").Value = 1066.7142
$ws.Range("M3").Value = -952.7141999999999
$ws.Range("H4").Value = 231.71428
$ws.Range("I4").Value = 278
$ws.Range("K4").Value = 278
$ws.Range("M4").Value = -163
$ws.Range("H22").Value = 1373.7142
$ws.Range("I22").Value = 1373.7142
$ws.Range("K22").Value = 1373.7142
$ws.Range("M22").Value = -1200.7142
$ws.Range("H24").Value = 2300
$ws.Range("I24").Value = 2300
$ws.Range("K24").Value = 2300
$ws.Range("M24").Value = -2065
$ws.Range("H134").Value = 3050.75
$ws.Range("I134").Value = 1839.3846
$ws.Range("K134").Value = 5518.1538
$ws.Range("M134").Value = -2983.1538

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1550
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H50").Value = 57645.625
$ws.Range("I50").Value = 38293.25
$ws.Range("J50").Value = 76998
$ws.Range("K50").Value = 38293.25
$ws.Range("L50").Value = 76998
$ws.Range("M50").Value = -37668.25
$ws.Range("N50").Value = -78248
$ws.Range("H132").Value = 2963.7856
$ws.Range("I132").Value = 2653.7693
$ws.Range("K132").Value = 7961.3079
$ws.Range("M132").Value = -5431.3079

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 19.2
$ws.Range("I7").Value = 16.428572
$ws.Range("J7").Value = 25.666666
$ws.Range("K7").Value = 49.28571599999999
$ws.Range("L7").Value = 76.99999800000001
$ws.Range("M7").Value = 62.71428400000001
$ws.Range("N7").Value = -300.999998
$ws.Range("H34").Value = 5417.1333
$ws.Range("I34").Value = 685.8
$ws.Range("J34").Value = 7782.8
$ws.Range("K34").Value = 2057.4
$ws.Range("L34").Value = 23348.4
$ws.Range("M34").Value = -1973.4
$ws.Range("N34").Value = -23516.4
$ws.Range("H39").Value = 10761.833
$ws.Range("J39").Value = 12155.6
$ws.Range("L39").Value = 36466.8
$ws.Range("N39").Value = -37054.8
$ws.Range("H55").Value = 4797.5
$ws.Range("I55").Value = 1333.3334
$ws.Range("J55").Value = 6876
$ws.Range("K55").Value = 4000.0002
$ws.Range("L55").Value = 20628
$ws.Range("M55").Value = -3823.0002
$ws.Range("N55").Value = -20982
$ws.Range("H107").Value = 465.83334
$ws.Range("I107").Value = 234
$ws.Range("J107").Value = 532.0714
$ws.Range("K107").Value = 702
$ws.Range("L107").Value = 1596.2142
$ws.Range("M107").Value = 1218
$ws.Range("N107").Value = -5436.2142

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 20778
$ws.Range("J95").Value = 20778
$ws.Range("L95").Value = 20778
$ws.Range("N95").Value = -26270
$ws.Range("H105").Value = 55988.6
$ws.Range("J105").Value = 55988.6
$ws.Range("L105").Value = 55988.6
$ws.Range("N105").Value = -62976.6
$ws.Range("H122").Value = 6172.125
$ws.Range("I122").Value = 1475.4
$ws.Range("J122").Value = 14000
$ws.Range("K122").Value = 4426.200000000001
$ws.Range("L122").Value = 42000
$ws.Range("M122").Value = -1976.200000000001
$ws.Range("N122").Value = -46900
$ws.Range("H126").Value = 4032.7273
$ws.Range("I126").Value = 3230.1428
$ws.Range("J126").Value = 5437.25
$ws.Range("K126").Value = 9690.428400000001
$ws.Range("L126").Value = 16311.75
$ws.Range("M126").Value = -7220.428400000001
$ws.Range("N126").Value = -21251.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1458.3529
$ws.Range("I22").Value = 690.3333
$ws.Range("K22").Value = 690.3333
$ws.Range("M22").Value = -395.3333
$ws.Range("H27").Value = 1458.3529
$ws.Range("I27").Value = 690.3333
$ws.Range("K27").Value = 690.3333
$ws.Range("M27").Value = -583.3333
$ws.Range("H46").Value = 3027.647
$ws.Range("I46").Value = 852.5
$ws.Range("K46").Value = 852.5
$ws.Range("M46").Value = -664.5
$ws.Range("H105").Value = 49999
$ws.Range("J105").Value = 49999
$ws.Range("L105").Value = 49999
$ws.Range("N105").Value = -56987
$ws.Range("H122").Value = 7580.727
$ws.Range("I122").Value = 3797.2
$ws.Range("J122").Value = 10733.667
$ws.Range("K122").Value = 11391.6
$ws.Range("L122").Value = 32201.001
$ws.Range("M122").Value = -8941.599999999999
$ws.Range("N122").Value = -37101.001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 70373.5
$ws.Range("J97").Value = 70373.5
$ws.Range("L97").Value = 70373.5
